# Update "想去人数" (number of people interested) counts for a few events
# across the workbook's sheets, per the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2046   # was 2041
$ws1.Range("F5").Value = 1130   # was 1126

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 16     # was 15

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2046   # was 2041
$ws4.Range("F5").Value = 16     # was 15
$ws4.Range("F7").Value = 1130   # was 1126
